# Fruta / hortaliza, semanal
# Re-order the weekly price rows (2-8) by updating the Date / Volumen /
# Precio minimo / Precio maximo / Precio promedio ponderado / Precio $/Kg
# columns (D, J, K, L, M, P) to reflect the new row order. All other
# columns stay identical across rows, so only these six columns need to
# be rewritten per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (row number => D, J, K, L, M, P)
$data = @{
    2 = @(44362, 120, 8000, 9000, 8500, 142)
    3 = @(44382, 160, 7000, 8000, 7438, 124)
    4 = @(44281, 120, 5500, 6000, 5750, 96)
    5 = @(44242, 160, 5000, 5500, 5250, 88)
    6 = @(44400, 120, 9000, 10000, 9500, 158)
    7 = @(44494, 120, 5000, 6000, 5500, 92)
    8 = @(44421, 100, 8000, 9000, 8500, 142)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
    $ws.Range("K$row").Value = $vals[2]
    $ws.Range("L$row").Value = $vals[3]
    $ws.Range("M$row").Value = $vals[4]
    $ws.Range("P$row").Value = $vals[5]
}
